$d = $word.ActiveDocument

# 1. Update the ThingSpeak channel number shown in the hyperlink display text
#    (the hyperlink's underlying relationship/address is left untouched,
#    matching the source edit which only retyped the visible text).
$d.Content.Find.Execute("channels/2232092", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "channels/2233511", 2) | Out-Null

# 2. Rename the "Assumptions" outline bullet
$d.Content.Find.Execute("Assumptions", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Connections/Subscriptions management", 2) | Out-Null

# 3. The report outline list ends with an empty bullet paragraph; fill it in
#    and append the remaining new outline bullets after it.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Text = "Retransmissions management"

$newItems = @(
    "Messages queue implementation",
    "Additional support variables and assumptions",
    "Debug channels"
)

$current = $d.Paragraphs.Item($d.Paragraphs.Count)
foreach ($item in $newItems) {
    $current.Range.InsertParagraphAfter()
    $current = $d.Paragraphs.Item($d.Paragraphs.Count)
    $current.Range.Text = $item
}
